$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4633537
$ws.Range("I51").Value = 3909.8462
$ws.Range("J51").Value = 16670567
$ws.Range("K51").Value = 3909.8462
$ws.Range("L51").Value = 16670567
$ws.Range("M51").Value = -3425.8462
$ws.Range("N51").Value = -16671535

$ws.Range("H74").Value = 5071.5713
$ws.Range("J74").Value = 5333.3335
$ws.Range("L74").Value = 5333.3335
$ws.Range("N74").Value = -7205.3335

$ws.Range("H77").Value = 5071.5713
$ws.Range("J77").Value = 5333.3335
$ws.Range("L77").Value = 26666.6675
$ws.Range("N77").Value = -36026.6675

$ws.Range("H95").Value = 34749.332
$ws.Range("J95").Value = 34749.332
$ws.Range("L95").Value = 34749.332
$ws.Range("N95").Value = -40241.332

$ws.Range("H113").Value = 95908.45
$ws.Range("J113").Value = 5665.8887
$ws.Range("L113").Value = 5665.8887
$ws.Range("N113").Value = -12173.8887

$ws.Range("H116").Value = 8469.046
$ws.Range("I116").Value = 10429.134
$ws.Range("J116").Value = 4268.857
$ws.Range("K116").Value = 10429.134
$ws.Range("L116").Value = 4268.857
$ws.Range("M116").Value = -6987.134
$ws.Range("N116").Value = -11152.857

$ws.Range("H132").Value = 2726.5151
$ws.Range("I132").Value = 1887.9259
$ws.Range("J132").Value = 6500.1665
$ws.Range("K132").Value = 5663.7777
$ws.Range("L132").Value = 19500.4995
$ws.Range("M132").Value = -3133.7777
$ws.Range("N132").Value = -24560.4995

$ws.Range("H137").Value = 2464.45
$ws.Range("J137").Value = 3750.75
$ws.Range("L137").Value = 11252.25
$ws.Range("N137").Value = -16352.25

$ws.Range("H141").Value = 5664.3125
$ws.Range("I141").Value = 3616.5
$ws.Range("K141").Value = 10849.5
$ws.Range("M141").Value = -5669.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5615.5
$ws.Range("I45").Value = 6975.048
$ws.Range("K45").Value = 6975.048
$ws.Range("M45").Value = -6598.048

$ws.Range("H61").Value = 2622.25
$ws.Range("I61").Value = 1931.625
$ws.Range("K61").Value = 1931.625
$ws.Range("M61").Value = -1719.625

$ws.Range("H74").Value = 3700.1082
$ws.Range("I74").Value = 3651.9678
$ws.Range("J74").Value = 3948.8333
$ws.Range("K74").Value = 3651.9678
$ws.Range("L74").Value = 3948.8333
$ws.Range("M74").Value = -2777.9678
$ws.Range("N74").Value = -5696.8333

$ws.Range("H77").Value = 3700.1082
$ws.Range("I77").Value = 3651.9678
$ws.Range("J77").Value = 3948.8333
$ws.Range("K77").Value = 18259.839
$ws.Range("L77").Value = 19744.1665
$ws.Range("M77").Value = -13891.839
$ws.Range("N77").Value = -28480.1665

$ws.Range("H132").Value = 2869.3
$ws.Range("I132").Value = 2442
$ws.Range("J132").Value = 3866.3333
$ws.Range("K132").Value = 7326
$ws.Range("L132").Value = 11598.9999
$ws.Range("M132").Value = -4796
$ws.Range("N132").Value = -16658.9999

$ws.Range("H136").Value = 2622.25
$ws.Range("I136").Value = 1931.625
$ws.Range("K136").Value = 5794.875
$ws.Range("M136").Value = -3244.875

$ws.Range("H138").Value = 106998.5
$ws.Range("J138").Value = 106998.5
$ws.Range("L138").Value = 106998.5
$ws.Range("N138").Value = -117278.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1334.4
$ws.Range("I134").Value = 1049.8823
$ws.Range("K134").Value = 3149.6469
$ws.Range("M134").Value = -614.6468999999997

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 45457412
$ws.Range("I7").Value = 66668760
$ws.Range("K7").Value = 66668760
$ws.Range("M7").Value = -66668647

$ws.Range("H28").Value = 28821.5
$ws.Range("J28").Value = 28821.5
$ws.Range("L28").Value = 28821.5
$ws.Range("N28").Value = -29311.5

$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()

$ws.Range("H58").Value = 1170.7368
$ws.Range("I58").Value = 1146.5
$ws.Range("J58").Value = 1300
$ws.Range("K58").Value = 1146.5
$ws.Range("L58").Value = 1300
$ws.Range("M58").Value = -943.5
$ws.Range("N58").Value = -1706

$ws.Range("H120").Value = 22898
$ws.Range("J120").Value = 22898
$ws.Range("L120").Value = 22898
$ws.Range("N120").Value = -30156

$ws.Range("H122").Value = 111912.89
$ws.Range("I122").Value = 125714.5
$ws.Range("K122").Value = 377143.5
$ws.Range("M122").Value = -374693.5

$ws.Range("H136").Value = 1170.7368
$ws.Range("I136").Value = 1146.5
$ws.Range("J136").Value = 1300
$ws.Range("K136").Value = 3439.5
$ws.Range("L136").Value = 3900
$ws.Range("M136").Value = -889.5
$ws.Range("N136").Value = -9000

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 2185262.8
$ws.Range("I32").Value = 350
$ws.Range("K32").Value = 1050
$ws.Range("M32").Value = -767

$ws.Range("H107").Value = 733.3333
$ws.Range("J107").Value = 200
$ws.Range("L107").Value = 600
$ws.Range("N107").Value = -4440

$ws.Range("H113").Value = 3000
$ws.Range("J113").Value = 3000
$ws.Range("L113").Value = 9000
$ws.Range("N113").Value = -13340

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3834.2727
$ws.Range("I113").Value = 5438.8
$ws.Range("K113").Value = 5438.8
$ws.Range("M113").Value = -3268.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 16543.39
$ws.Range("I7").Value = 21041
$ws.Range("K7").Value = 21041
$ws.Range("M7").Value = -20929

$ws.Range("H22").Value = 4329.75
$ws.Range("I22").Value = 1563
$ws.Range("J22").Value = 5989.8
$ws.Range("K22").Value = 1563
$ws.Range("L22").Value = 5989.8
$ws.Range("M22").Value = -1268
$ws.Range("N22").Value = -6579.8

$ws.Range("H27").Value = 4329.75
$ws.Range("I27").Value = 1563
$ws.Range("J27").Value = 5989.8
$ws.Range("K27").Value = 1563
$ws.Range("L27").Value = 5989.8
$ws.Range("M27").Value = -1456
$ws.Range("N27").Value = -6203.8

$ws.Range("H46").Value = 3386.625
$ws.Range("I46").Value = 2918.2
$ws.Range("J46").Value = 4167.3335
$ws.Range("K46").Value = 2918.2
$ws.Range("L46").Value = 4167.3335
$ws.Range("M46").Value = -2730.2
$ws.Range("N46").Value = -4543.3335

$ws.Range("H61").Value = 500949.5
$ws.Range("I61").Value = 1000000
$ws.Range("K61").Value = 1000000
$ws.Range("M61").Value = -999798

$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()

$ws.Range("H113").Value = 500949.5
$ws.Range("I113").Value = 1000000
$ws.Range("K113").Value = 1000000
$ws.Range("M113").Value = -997830

$ws.Range("H122").Value = 103991.55
$ws.Range("I122").Value = 114768.39
$ws.Range("K122").Value = 344305.17
$ws.Range("M122").Value = -341855.17

$ws.Range("H126").Value = 16543.39
$ws.Range("I126").Value = 21041
$ws.Range("K126").Value = 63123
$ws.Range("M126").Value = -60653

$ws.Range("H134").Value = 74666
$ws.Range("J134").Value = 74666
$ws.Range("L134").Value = 74666
$ws.Range("N134").Value = -84806

$ws.Range("H136").Value = 5069.533
$ws.Range("J136").Value = 6879
$ws.Range("L136").Value = 20637
$ws.Range("N136").Value = -25737

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 13159008
$ws.Range("I107").Value = 1114.36
$ws.Range("K107").Value = 3343.08
$ws.Range("M107").Value = -1423.08

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

$ws.Range("H126").Value = 2043.6
$ws.Range("I126").Value = 1739.3334
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 5218.0002
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -2748.0002
$ws.Range("N126").Value = -12440

$ws.Range("H136").Value = 2045.1578
$ws.Range("I136").Value = 1077.6364
$ws.Range("J136").Value = 3375.5
$ws.Range("K136").Value = 3232.9092
$ws.Range("L136").Value = 10126.5
$ws.Range("M136").Value = -682.9092000000001
$ws.Range("N136").Value = -15226.5

